# Update the "dSF" column (F) with recalculated values after repulling data.
# Rows correspond to the worksheet rows: row 1 = header, rows 2-29 = data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -6
    5  = -3
    6  = 2
    7  = -1
    8  = 1
    9  = 6
    10 = 1
    11 = -2
    12 = 3
    13 = 4
    14 = 2
    15 = 2
    16 = 3
    17 = 5
    18 = -3
    19 = -2
    20 = -2
    21 = 2
    23 = 1
    24 = 2
    25 = 4
    26 = -3
    27 = 2
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
